$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data to append at row 26 (weekly update of fruit/vegetable prices)
$row = 26

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"

# Column D: date value, formatted as date (serial 45191 -> 2023-09-22)
# Set number format before assigning the value so Excel reuses the existing
# date style instead of creating a new (unused) default date style.
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 4).Value = (Get-Date -Year 2023 -Month 9 -Day 22 -Hour 0 -Minute 0 -Second 0)

$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 100112039
$ws.Cells.Item($row, 7).Value = "Ciboulette"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 100
$ws.Cells.Item($row, 11).Value = 2500
$ws.Cells.Item($row, 12).Value = 2500
$ws.Cells.Item($row, 13).Value = 2500
$ws.Cells.Item($row, 14).Value = "`$/docena de atados"
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 833
$ws.Cells.Item($row, 17).Value = 3
$ws.Cells.Item($row, 18).Value = "Hortaliza"
